$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns with latest symbol data.
# Values are written with a leading apostrophe to force them to remain as literal
# text (matching the original inlineStr text cells) rather than being auto-converted
# to numeric/percentage values by Excel.

$ws.Range("D2").Value = "'328.93"
$ws.Range("E2").Value = "'0.15%"
$ws.Range("D3").Value = "'44.41"
$ws.Range("E3").Value = "'0.01%"
$ws.Range("D4").Value = "'5.494"
$ws.Range("E4").Value = "'-1.77%"
$ws.Range("D5").Value = "'0.08072"
$ws.Range("E5").Value = "'-0.19%"
$ws.Range("D6").Value = "'2.056"
$ws.Range("E6").Value = "'-0.07%"
$ws.Range("D7").Value = "'0.9590"
$ws.Range("E7").Value = "'0.64%"
$ws.Range("E8").Value = "'-3.43%"
$ws.Range("E9").Value = "'1.39%"
$ws.Range("D10").Value = "'10.26"
$ws.Range("E10").Value = "'-0.09%"
$ws.Range("D11").Value = "'0.09925"
$ws.Range("E11").Value = "'1.24%"
$ws.Range("D12").Value = "'0.04713"
$ws.Range("E12").Value = "'4.20%"
$ws.Range("D13").Value = "'0.1053"
$ws.Range("E13").Value = "'-1.35%"
$ws.Range("D14").Value = "'0.001262"
$ws.Range("E14").Value = "'-1.66%"
$ws.Range("D15").Value = "'0.04117"
$ws.Range("E15").Value = "'-1.94%"
$ws.Range("D16").Value = "'0.006125"
$ws.Range("E16").Value = "'4.00%"
$ws.Range("E17").Value = "'-0.88%"
$ws.Range("D18").Value = "'4.435"
$ws.Range("E18").Value = "'2.90%"
$ws.Range("D19").Value = "'2.621"
$ws.Range("E19").Value = "'4.11%"
$ws.Range("D20").Value = "'0.3316"
$ws.Range("E20").Value = "'-4.75%"
$ws.Range("D21").Value = "'0.1393"
$ws.Range("E21").Value = "'-1.06%"
$ws.Range("D22").Value = "'0.2581"
$ws.Range("E22").Value = "'2.98%"
$ws.Range("D23").Value = "'0.001312"
$ws.Range("E23").Value = "'5.25%"
$ws.Range("D24").Value = "'0.004351"
$ws.Range("E24").Value = "'0.15%"
$ws.Range("D25").Value = "'0.0001284"
$ws.Range("E25").Value = "'7.67%"
$ws.Range("D26").Value = "'0.0003749"
$ws.Range("E26").Value = "'-5.79%"
$ws.Range("D38").Value = "'0.02629"
$ws.Range("E38").Value = "'-1.42%"
$ws.Range("D39").Value = "'0.05620"
$ws.Range("E39").Value = "'0.99%"
$ws.Range("D40").Value = "'0.007596"
$ws.Range("E40").Value = "'0.35%"
$ws.Range("D41").Value = "'0.1400"
$ws.Range("E41").Value = "'-0.63%"
$ws.Range("D42").Value = "'0.007403"
$ws.Range("E42").Value = "'-6.94%"
$ws.Range("D43").Value = "'0.001912"
$ws.Range("E43").Value = "'-5.38%"
$ws.Range("D44").Value = "'0.008738"
$ws.Range("E44").Value = "'3.87%"
$ws.Range("D45").Value = "'0.00007118"
$ws.Range("E45").Value = "'-0.68%"
$ws.Range("D46").Value = "'0.00000000752"
$ws.Range("E46").Value = "'0.01%"
$ws.Range("D47").Value = "'0.0005818"
$ws.Range("E47").Value = "'0.11%"
$ws.Range("D48").Value = "'0.002526"
$ws.Range("E48").Value = "'11.23%"
$ws.Range("D49").Value = "'0.003522"
$ws.Range("E49").Value = "'-14.73%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.01%"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.01%"
